# Apply the CASPER_1.0_Manual.docx edit described by the diff:
#   1. The first numbered-list item ("Placeholder") becomes real content:
#      "Argument 1 is always the executable file name." and the
#      auto-managed "_GoBack" bookmark now sits right after that run.
#   2. The "_GoBack" bookmark that used to sit between " OffTarget" and
#      ".py" (an artifact of where the author's cursor last was) is
#      removed, since Word only ever keeps one "_GoBack" bookmark at a
#      time - it simply moved to the new edit location above.

$d = $word.ActiveDocument

# --- Step 1: turn the placeholder bullet into real text -------------------
$d.Content.Find.Execute("Placeholder", $true, $false, $false, $false, $false, $true, 1, $false, "Argument 1 is always the executable file name.", 2)

# --- Step 2: drop the old "_GoBack" bookmark near " OffTarget.py" ---------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3: re-create "_GoBack" right after the new sentence -------------
# Locate the end of the sentence we just inserted.
$found = $d.Content
$found.Find.Execute("Argument 1 is always the executable file name.")
$endPos = $found.End

# A collapsed Range sitting exactly at a paragraph's end (i.e. right before
# its paragraph mark) confuses this host's Bookmarks.Add, so we nudge it:
# insert a throwaway character after the target spot, add the bookmark
# while the Range is safely mid-paragraph, then delete the throwaway
# character again. The bookmark (zero-width) stays put as its neighboring
# text is removed.
$guard = $d.Range($endPos, $endPos)
$guard.InsertAfter("X")

$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$d.Range($endPos, $endPos + 1).Delete()
